$d = $word.ActiveDocument

$replacements = @(
    @("54×42=", "30×52="),
    @("28×91=", "76×82="),
    @("74×52=", "27×84="),
    @("52×36=", "63×27="),
    @("97×79=", "21×54="),
    @("83×50=", "39×44="),
    @("42×46=", "16×86="),
    @("98×80=", "98×72="),
    @("85×86=", "33×14="),
    @("63×59=", "61×30="),
    @("82×22=", "38×38="),
    @("76×90=", "67×71="),
    @("28×53=", "56×22="),
    @("24×56=", "77×80="),
    @("47×22=", "89×60="),
    @("16×98=", "27×21="),
    @("51×32=", "80×68="),
    @("61×84=", "67×70="),
    @("45×60=", "37×27="),
    @("86×20=", "42×75="),
    @("40×46=", "75×54="),
    @("51×29=", "19×65="),
    @("42×92=", "65×17="),
    @("91×57=", "39×16="),
    @("47×68=", "85×35=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
